$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Remove the duplicated "Contact" / "No display for ContactDetail" row (row 11).
# Deleting it shifts rows 12-22 up to 11-21 while keeping their existing cell
# styles untouched (no new style gets synthesized), which lines the remaining
# rows up with the target layout (dimension becomes A1:B21).
$ws.Rows.Item(11).Delete()

# Publisher value: was empty -> "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# The row that used to hold the (now removed) second "Contact" row is reused
# for the new "Jurisdiction" / "United States of America" property.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Case Sensitive value: was empty -> "true" (stored as literal text, not an
# Excel boolean). Assigning the bare word "true" via Value/Formula gets
# auto-coerced into a boolean TRUE cell, and forcing text with a leading
# apostrophe marks the cell with a "quote prefix" style. Instead, compute
# the text "true" as a formula result in a scratch cell (a formula's string
# result is plain text, not a boolean) and paste only its value into place;
# this keeps B14's original style/border/wrap formatting untouched.
$ws.Range("Z1").Formula = "=""true"""
$ws.Range("Z1").Copy()
$ws.Range("B14").PasteSpecial(-4163)  # xlPasteValues
$ws.Columns.Item(26).Delete()
$excel.CutCopyMode = $false
